$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.514.36'
$ws.Range('E2').Value = '  -3.02%  '
$ws.Range('D3').Value = '3.269.59'
$ws.Range('E3').Value = '  -5.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.12'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.89'
$ws.Range('E6').Value = '  -9.70%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.259.53'
$ws.Range('E8').Value = '  -5.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  -8.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  -10.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.76'
$ws.Range('E11').Value = '  -4.70%  '
$ws.Range('E12').Value = '  -9.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.67'
$ws.Range('E13').Value = '  -12.74%  '
$ws.Range('E14').Value = '  -7.93%  '
$ws.Range('D15').Value = '3.790.38'
$ws.Range('E15').Value = '  -5.62%  '
$ws.Range('D16').Value = '67.503.49'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').Value = '3.266.19'
$ws.Range('E17').Value = '  -5.85%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.114'
$ws.Range('E18').Value = '  -5.10%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '533.88'
$ws.Range('E19').Value = '  -8.64%  '
$ws.Range('E20').Value = '  -12.32%  '
$ws.Range('E21').Value = '  -12.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.762'
$ws.Range('E22').Value = '  -10.52%  '
$ws.Range('E23').Value = '  -11.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.79'
$ws.Range('E24').Value = '  -9.84%  '
$ws.Range('E25').Value = '  -10.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -9.75%  '
$ws.Range('E28').Value = '  -5.96%  '
$ws.Range('E29').Value = '  -11.86%  '
$ws.Range('E30').Value = '  -10.89%  '
$ws.Range('E31').Value = '  -4.03%  '
$ws.Range('E32').Value = '  -6.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.65'
$ws.Range('E33').Value = '  -14.90%  '
$ws.Range('E34').Value = '  -12.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '523.45'
$ws.Range('E35').Value = '  -9.02%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0455'
$ws.Range('E37').Value = '  -5.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.57'
$ws.Range('E38').Value = '  -4.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0861'
$ws.Range('E39').Value = '  -10.20%  '
$ws.Range('E40').Value = '  -14.70%  '
$ws.Range('E41').Value = '  -9.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.82'
$ws.Range('E42').Value = '  -9.82%  '
$ws.Range('D43').Value = '2.942.95'
$ws.Range('E43').Value = '  -9.41%  '
$ws.Range('E44').Value = '  -9.01%  '
$ws.Range('E45').Value = '  -14.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.20'
$ws.Range('E46').Value = '  -7.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.83'
$ws.Range('E47').Value = '  -12.73%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('E49').Value = '  -14.87%  '
$ws.Range('E50').Value = '  -9.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '123.39'
$ws.Range('E51').Value = '  -7.53%  '
